$wb = $excel.ActiveWorkbook

# Sheet 1: LP1912 - add a new data row and update header info
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 01:55:40"
$ws1.Range("A3").Value = "Total filas: 4"

$ws1.Cells.Item(9, 1).Value = "01:55:40"
$ws1.Cells.Item(9, 2).Value = "03:48"
$ws1.Cells.Item(9, 3).Value = "14_ABASTO"
$ws1.Cells.Item(9, 4).Value = 113
$ws1.Cells.Item(9, 5).Value = "LP1912"

# Sheet 2: LP1912-215 - only update the "last updated" timestamp
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 01:55:40"

# Sheet 3: 6203-6173 - only update the "last updated" timestamp
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 01:55:40"
